$wb = $excel.ActiveWorkbook
$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

$schedule.Range("E4").Value = 460.8246629999999
$schedule.Range("F4").Value = 30.47782162698412
$schedule.Range("E5").Value = -38.58871575000001
$schedule.Range("F5").Value = -1.134294995590829
$detailed.Range("B35").Value = -2.47052
$detailed.Range("B36").Value = 0
$detailed.Range("B37").Value = -2.58035
$detailed.Range("B38").Value = 3.34464
$detailed.Range("C38").Value = "historical"
$detailed.Range("B39").Value = 38.55471
$detailed.Range("B40").Value = 57.18142
$detailed.Range("B41").Value = 64.53386
$detailed.Range("B42").Value = 65
$detailed.Range("B44").Value = 65.37398
$detailed.Range("B45").Value = 65.74816
$detailed.Range("B46").Value = 59.50728
$detailed.Range("B47").Value = 64.8901
$detailed.Range("B48").Value = 59.04126
$detailed.Range("B49").Value = 63.32463
$detailed.Range("B50").Value = 61.17037
$detailed.Range("B51").Value = 60.8643
$detailed.Range("B57").Value = 61.1294
$detailed.Range("B58").Value = 63.69192
$detailed.Range("B59").Value = 65
$detailed.Range("B60").Value = 65.21083
$detailed.Range("B61").Value = 75.57516
$detailed.Range("B62").Value = 77.94
$detailed.Range("B63").Value = 75.76638
$detailed.Range("B64").Value = 56.98
$detailed.Range("B65").Value = 8.85098
$detailed.Range("B67").Value = 0.01108
$detailed.Range("B68").Value = -2.57556
$detailed.Range("B69").Value = -6
$detailed.Range("B70").Value = -7.01255
$detailed.Range("B71").Value = -6.7961
$detailed.Range("B72").Value = -8.13832
$detailed.Range("B73").Value = -7.78537
$detailed.Range("B74").Value = -7.78335
$detailed.Range("B75").Value = -10.7155
$detailed.Range("B77").Value = -8.28947
$detailed.Range("B78").Value = -8.179539999999999
$detailed.Range("B79").Value = -8.558400000000001
$detailed.Range("B80").Value = -7.86778
$detailed.Range("B81").Value = -6.42829
$detailed.Range("B82").Value = -2.80319
$detailed.Range("B83").Value = -5.51
$detailed.Range("B84").Value = -1.01549
$detailed.Range("B85").Value = 9.46617
$detailed.Range("B86").Value = 9.59477
$detailed.Range("B87").Value = 31.01719
$detailed.Range("B88").Value = 57.3
$detailed.Range("B90").Value = 60.56861
$detailed.Range("B91").Value = 58.94198
$detailed.Range("B92").Value = 58.58436
$detailed.Range("B93").Value = 57.84212
$detailed.Range("B94").Value = 57.3
$detailed.Range("B95").Value = 58.86227
$detailed.Range("B96").Value = 61.45254
